$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.54348886013031
$ws.Range("B1").Value = 2.764216423034668
$ws.Range("C1").Value = 6.390307903289795
$ws.Range("D1").Value = 1.588533401489258
$ws.Range("E1").Value = 0.8652488589286804
